$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2114093959731544
$ws.Range("C2").Value = 0.5100671140939598
$ws.Range("J2").Value = 0.01006711409395973
$ws.Range("P2").Value = 0.1711409395973154
$ws.Range("S2").Value = 0.09731543624161074
$ws.Range("B3").Value = 0.01204819277108434
$ws.Range("C3").Value = 0.03012048192771084
$ws.Range("J3").Value = 0.02409638554216868
$ws.Range("P3").Value = 0.7530120481927711
$ws.Range("S3").Value = 0.1807228915662651
$ws.Range("P4").Value = 0.65
$ws.Range("S4").Value = 0.35
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.02857142857142857
$ws.Range("D6").Value = 0.01224489795918367
$ws.Range("E6").Value = 0.004081632653061225
$ws.Range("F6").Value = 0.06122448979591837
$ws.Range("J6").Value = 0.2408163265306122
$ws.Range("O6").Value = 0.00816326530612245
$ws.Range("Q6").Value = 0.1918367346938775
$ws.Range("R6").Value = 0.09387755102040816
$ws.Range("S6").Value = 0.3591836734693877
$ws.Range("B7").Value = 0.09523809523809523
$ws.Range("D7").Value = 0.02380952380952381
$ws.Range("E7").Value = 0.005952380952380952
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.119047619047619
$ws.Range("O7").Value = 0.005952380952380952
$ws.Range("Q7").Value = 0.1785714285714286
$ws.Range("R7").Value = 0.07738095238095238
$ws.Range("S7").Value = 0.4464285714285715
$ws.Range("B8").Value = 0.09188034188034189
$ws.Range("D8").Value = 0.01923076923076923
$ws.Range("E8").Value = 0.002136752136752137
$ws.Range("F8").Value = 0.05982905982905983
$ws.Range("J8").Value = 0.0811965811965812
$ws.Range("O8").Value = 0.01068376068376068
$ws.Range("Q8").Value = 0.1858974358974359
$ws.Range("R8").Value = 0.08547008547008547
$ws.Range("S8").Value = 0.4636752136752137
$ws.Range("B9").Value = 0.09259259259259259
$ws.Range("D9").Value = 0.02592592592592593
$ws.Range("E9").Value = 0.003703703703703704
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.07407407407407407
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.06296296296296296
$ws.Range("S9").Value = 0.4851851851851852
$ws.Range("B10").Value = 0.1026022304832714
$ws.Range("D10").Value = 0.01338289962825279
$ws.Range("E10").Value = 0.0007434944237918215
$ws.Range("F10").Value = 0.07211895910780669
$ws.Range("J10").Value = 0.09814126394052045
$ws.Range("O10").Value = 0.02007434944237918
$ws.Range("Q10").Value = 0.2423791821561338
$ws.Range("R10").Value = 0.07732342007434945
$ws.Range("S10").Value = 0.3732342007434944
$ws.Range("F11").Value = 0.003703703703703704
$ws.Range("G11").Value = 0.1185185185185185
$ws.Range("J11").Value = 0.1148148148148148
$ws.Range("K11").Value = 0.1888888888888889
$ws.Range("L11").Value = 0.5666666666666667
$ws.Range("S11").Value = 0.007407407407407408
$ws.Range("G12").Value = 0.7423312883435583
$ws.Range("J12").Value = 0.1595092024539877
$ws.Range("K12").Value = 0.0245398773006135
$ws.Range("L12").Value = 0.05521472392638037
$ws.Range("S12").Value = 0.01840490797546012
$ws.Range("G13").Value = 0.8387096774193549
$ws.Range("J13").Value = 0.1290322580645161
$ws.Range("S13").Value = 0.03225806451612903
$ws.Range("F15").Value = 0.02164502164502164
$ws.Range("H15").Value = 0.1601731601731602
$ws.Range("I15").Value = 0.08658008658008658
$ws.Range("J15").Value = 0.3593073593073593
$ws.Range("K15").Value = 0.04761904761904762
$ws.Range("M15").Value = 0.008658008658008658
$ws.Range("O15").Value = 0.04761904761904762
$ws.Range("S15").Value = 0.2683982683982684
$ws.Range("F16").Value = 0.03571428571428571
$ws.Range("H16").Value = 0.1479591836734694
$ws.Range("I16").Value = 0.08673469387755102
$ws.Range("J16").Value = 0.4744897959183674
$ws.Range("K16").Value = 0.08673469387755102
$ws.Range("M16").Value = 0.02040816326530612
$ws.Range("N16").Value = 0.00510204081632653
$ws.Range("O16").Value = 0.05612244897959184
$ws.Range("S16").Value = 0.08673469387755102
$ws.Range("F17").Value = 0.02209944751381215
$ws.Range("H17").Value = 0.1657458563535912
$ws.Range("I17").Value = 0.1141804788213628
$ws.Range("J17").Value = 0.4475138121546962
$ws.Range("K17").Value = 0.06998158379373849
$ws.Range("M17").Value = 0.009208103130755065
$ws.Range("O17").Value = 0.07366482504604052
$ws.Range("S17").Value = 0.09760589318600368
$ws.Range("F18").Value = 0.04123711340206185
$ws.Range("H18").Value = 0.1391752577319588
$ws.Range("I18").Value = 0.1185567010309278
$ws.Range("J18").Value = 0.4175257731958763
$ws.Range("K18").Value = 0.1134020618556701
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("N18").Value = 0.005154639175257732
$ws.Range("O18").Value = 0.05670103092783505
$ws.Range("S18").Value = 0.08762886597938144
$ws.Range("F19").Value = 0.01466275659824047
$ws.Range("H19").Value = 0.2111436950146628
$ws.Range("I19").Value = 0.1085043988269795
$ws.Range("J19").Value = 0.3870967741935484
$ws.Range("K19").Value = 0.09017595307917889
$ws.Range("M19").Value = 0.01319648093841642
$ws.Range("N19").Value = 0.0007331378299120235
$ws.Range("O19").Value = 0.07038123167155426
$ws.Range("S19").Value = 0.1041055718475073
